$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Natmi following Dr Hou advice
# Update ligand/receptor expressing-cell counts (E,K: 1 -> 3) and recompute
# dependent expression/specificity metrics for rows 2-17 of the
# Ecm1-Itgb4 LR-pair sheet.

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 13.671071
$ws.Cells.Item(2, 8).Value = 41.013213
$ws.Cells.Item(2, 9).Value = 0.09691937964924315
$ws.Cells.Item(2, 10).Value = 0.09691937964924313
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 5.378434333333334
$ws.Cells.Item(2, 14).Value = 16.135303
$ws.Cells.Item(2, 15).Value = 0.5423901270513668
$ws.Cells.Item(2, 16).Value = 0.5423901270513669
$ws.Cells.Item(2, 17).Value = 73.52895763983767
$ws.Cells.Item(2, 18).Value = 661.7606187585391
$ws.Cells.Item(2, 19).Value = 0.05256811464169264
$ws.Cells.Item(2, 20).Value = 0.05256811464169265
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 13.671071
$ws.Cells.Item(3, 8).Value = 41.013213
$ws.Cells.Item(3, 9).Value = 0.09691937964924315
$ws.Cells.Item(3, 10).Value = 0.09691937964924313
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 2.118986
$ws.Cells.Item(3, 14).Value = 6.356958
$ws.Cells.Item(3, 15).Value = 0.213689898310568
$ws.Cells.Item(3, 16).Value = 0.213689898310568
$ws.Cells.Item(3, 17).Value = 28.968808054006
$ws.Cells.Item(3, 18).Value = 260.719272486054
$ws.Cells.Item(3, 19).Value = 0.0207106923815701
$ws.Cells.Item(3, 20).Value = 0.0207106923815701
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 13.671071
$ws.Cells.Item(4, 8).Value = 41.013213
$ws.Cells.Item(4, 9).Value = 0.09691937964924315
$ws.Cells.Item(4, 10).Value = 0.09691937964924313
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 0.3858953333333333
$ws.Cells.Item(4, 14).Value = 1.157686
$ws.Cells.Item(4, 15).Value = 0.0389157524110696
$ws.Cells.Item(4, 16).Value = 0.03891575241106961
$ws.Cells.Item(4, 17).Value = 5.275602500568667
$ws.Cells.Item(4, 18).Value = 47.480422505118
$ws.Cells.Item(4, 19).Value = 0.003771690582264404
$ws.Cells.Item(4, 20).Value = 0.003771690582264404
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 13.671071
$ws.Cells.Item(5, 8).Value = 41.013213
$ws.Cells.Item(5, 9).Value = 0.09691937964924315
$ws.Cells.Item(5, 10).Value = 0.09691937964924313
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 2.032857333333333
$ws.Cells.Item(5, 14).Value = 6.098572
$ws.Cells.Item(5, 15).Value = 0.2050042222269955
$ws.Cells.Item(5, 16).Value = 0.2050042222269956
$ws.Cells.Item(5, 17).Value = 27.79133693687066
$ws.Cells.Item(5, 18).Value = 250.122032431836
$ws.Cells.Item(5, 19).Value = 0.01986888204371599
$ws.Cells.Item(5, 20).Value = 0.01986888204371599
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 92.04504633333333
$ws.Cells.Item(6, 8).Value = 276.135139
$ws.Cells.Item(6, 9).Value = 0.6525420569034064
$ws.Cells.Item(6, 10).Value = 0.6525420569034063
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 5.378434333333334
$ws.Cells.Item(6, 14).Value = 16.135303
$ws.Cells.Item(6, 15).Value = 0.5423901270513668
$ws.Cells.Item(6, 16).Value = 0.5423901270513669
$ws.Cells.Item(6, 17).Value = 495.0582374124575
$ws.Cells.Item(6, 18).Value = 4455.524136712117
$ws.Cells.Item(6, 19).Value = 0.3539323691501988
$ws.Cells.Item(6, 20).Value = 0.3539323691501988
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 92.04504633333333
$ws.Cells.Item(7, 8).Value = 276.135139
$ws.Cells.Item(7, 9).Value = 0.6525420569034064
$ws.Cells.Item(7, 10).Value = 0.6525420569034063
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 2.118986
$ws.Cells.Item(7, 14).Value = 6.356958
$ws.Cells.Item(7, 15).Value = 0.213689898310568
$ws.Cells.Item(7, 16).Value = 0.213689898310568
$ws.Cells.Item(7, 17).Value = 195.0421645496847
$ws.Cells.Item(7, 18).Value = 1755.379480947162
$ws.Cells.Item(7, 19).Value = 0.1394416457830578
$ws.Cells.Item(7, 20).Value = 0.1394416457830578
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 92.04504633333333
$ws.Cells.Item(8, 8).Value = 276.135139
$ws.Cells.Item(8, 9).Value = 0.6525420569034064
$ws.Cells.Item(8, 10).Value = 0.6525420569034063
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 0.3858953333333333
$ws.Cells.Item(8, 14).Value = 1.157686
$ws.Cells.Item(8, 15).Value = 0.0389157524110696
$ws.Cells.Item(8, 16).Value = 0.03891575241106961
$ws.Cells.Item(8, 17).Value = 35.51975383648377
$ws.Cells.Item(8, 18).Value = 319.677784528354
$ws.Cells.Item(8, 19).Value = 0.02539416512426305
$ws.Cells.Item(8, 20).Value = 0.02539416512426305
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 92.04504633333333
$ws.Cells.Item(9, 8).Value = 276.135139
$ws.Cells.Item(9, 9).Value = 0.6525420569034064
$ws.Cells.Item(9, 10).Value = 0.6525420569034063
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 2.032857333333333
$ws.Cells.Item(9, 14).Value = 6.098572
$ws.Cells.Item(9, 15).Value = 0.2050042222269955
$ws.Cells.Item(9, 16).Value = 0.2050042222269956
$ws.Cells.Item(9, 17).Value = 187.1144474357231
$ws.Cells.Item(9, 18).Value = 1684.030026921508
$ws.Cells.Item(9, 19).Value = 0.1337738768458867
$ws.Cells.Item(9, 20).Value = 0.1337738768458867
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 25.57700433333333
$ws.Cells.Item(10, 8).Value = 76.731013
$ws.Cells.Item(10, 9).Value = 0.1813250324917975
$ws.Cells.Item(10, 10).Value = 0.1813250324917975
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 5.378434333333334
$ws.Cells.Item(10, 14).Value = 16.135303
$ws.Cells.Item(10, 15).Value = 0.5423901270513668
$ws.Cells.Item(10, 16).Value = 0.5423901270513669
$ws.Cells.Item(10, 17).Value = 137.5642382502155
$ws.Cells.Item(10, 18).Value = 1238.078144251939
$ws.Cells.Item(10, 19).Value = 0.09834890741081927
$ws.Cells.Item(10, 20).Value = 0.09834890741081927
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 25.57700433333333
$ws.Cells.Item(11, 8).Value = 76.731013
$ws.Cells.Item(11, 9).Value = 0.1813250324917975
$ws.Cells.Item(11, 10).Value = 0.1813250324917975
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 2.118986
$ws.Cells.Item(11, 14).Value = 6.356958
$ws.Cells.Item(11, 15).Value = 0.213689898310568
$ws.Cells.Item(11, 16).Value = 0.213689898310568
$ws.Cells.Item(11, 17).Value = 54.19731410427267
$ws.Cells.Item(11, 18).Value = 487.775826938454
$ws.Cells.Item(11, 19).Value = 0.03874732775433264
$ws.Cells.Item(11, 20).Value = 0.03874732775433264
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 25.57700433333333
$ws.Cells.Item(12, 8).Value = 76.731013
$ws.Cells.Item(12, 9).Value = 0.1813250324917975
$ws.Cells.Item(12, 10).Value = 0.1813250324917975
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 0.3858953333333333
$ws.Cells.Item(12, 14).Value = 1.157686
$ws.Cells.Item(12, 15).Value = 0.0389157524110696
$ws.Cells.Item(12, 16).Value = 0.03891575241106961
$ws.Cells.Item(12, 17).Value = 9.870046612879777
$ws.Cells.Item(12, 18).Value = 88.83041951591801
$ws.Cells.Item(12, 19).Value = 0.007056400070379942
$ws.Cells.Item(12, 20).Value = 0.007056400070379943
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 25.57700433333333
$ws.Cells.Item(13, 8).Value = 76.731013
$ws.Cells.Item(13, 9).Value = 0.1813250324917975
$ws.Cells.Item(13, 10).Value = 0.1813250324917975
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 2.032857333333333
$ws.Cells.Item(13, 14).Value = 6.098572
$ws.Cells.Item(13, 15).Value = 0.2050042222269955
$ws.Cells.Item(13, 16).Value = 0.2050042222269956
$ws.Cells.Item(13, 17).Value = 51.99440082371511
$ws.Cells.Item(13, 18).Value = 467.949607413436
$ws.Cells.Item(13, 19).Value = 0.03717239725626564
$ws.Cells.Item(13, 20).Value = 0.03717239725626564
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 9.762991666666666
$ws.Cells.Item(14, 8).Value = 29.288975
$ws.Cells.Item(14, 9).Value = 0.06921353095555306
$ws.Cells.Item(14, 10).Value = 0.06921353095555306
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 5.378434333333334
$ws.Cells.Item(14, 14).Value = 16.135303
$ws.Cells.Item(14, 15).Value = 0.5423901270513668
$ws.Cells.Item(14, 16).Value = 0.5423901270513669
$ws.Cells.Item(14, 17).Value = 52.50960957604723
$ws.Cells.Item(14, 18).Value = 472.586486184425
$ws.Cells.Item(14, 19).Value = 0.03754073584865614
$ws.Cells.Item(14, 20).Value = 0.03754073584865614
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 9.762991666666666
$ws.Cells.Item(15, 8).Value = 29.288975
$ws.Cells.Item(15, 9).Value = 0.06921353095555306
$ws.Cells.Item(15, 10).Value = 0.06921353095555306
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 2.118986
$ws.Cells.Item(15, 14).Value = 6.356958
$ws.Cells.Item(15, 15).Value = 0.213689898310568
$ws.Cells.Item(15, 16).Value = 0.213689898310568
$ws.Cells.Item(15, 17).Value = 20.68764265978333
$ws.Cells.Item(15, 18).Value = 186.18878393805
$ws.Cells.Item(15, 19).Value = 0.01479023239160748
$ws.Cells.Item(15, 20).Value = 0.01479023239160748
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 9.762991666666666
$ws.Cells.Item(16, 8).Value = 29.288975
$ws.Cells.Item(16, 9).Value = 0.06921353095555306
$ws.Cells.Item(16, 10).Value = 0.06921353095555306
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 0.3858953333333333
$ws.Cells.Item(16, 14).Value = 1.157686
$ws.Cells.Item(16, 15).Value = 0.0389157524110696
$ws.Cells.Item(16, 16).Value = 0.03891575241106961
$ws.Cells.Item(16, 17).Value = 3.767492923538889
$ws.Cells.Item(16, 18).Value = 33.90743631185
$ws.Cells.Item(16, 19).Value = 0.002693496634162205
$ws.Cells.Item(16, 20).Value = 0.002693496634162205
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 9.762991666666666
$ws.Cells.Item(17, 8).Value = 29.288975
$ws.Cells.Item(17, 9).Value = 0.06921353095555306
$ws.Cells.Item(17, 10).Value = 0.06921353095555306
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 2.032857333333333
$ws.Cells.Item(17, 14).Value = 6.098572
$ws.Cells.Item(17, 15).Value = 0.2050042222269955
$ws.Cells.Item(17, 16).Value = 0.2050042222269956
$ws.Cells.Item(17, 17).Value = 19.84676920485555
$ws.Cells.Item(17, 18).Value = 178.6209228437
$ws.Cells.Item(17, 19).Value = 0.01418906608112724
$ws.Cells.Item(17, 20).Value = 0.01418906608112724
